$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 42605.888182870367
$ws.Range("B3").Value = -30
$ws.Range("C3").Value = 55
$ws.Range("D3").Value = 42
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 86
$ws.Range("G3").Value = 23677
$ws.Range("H3").Value = 8871
$ws.Range("I3").Value = 1382
$ws.Range("J3").Value = 146
$ws.Range("K3").Value = 111
$ws.Range("L3").Value = 5
$ws.Range("M3").Value = 32
$ws.Range("N3").Value = "Noun"

$ws.Range("A3").NumberFormat = "m/d/yy h:mm"
